$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 882.125
$ws.Range("J41").Value = 983.0909
$ws.Range("L41").Value = 983.0909
$ws.Range("N41").Value = -1863.0909
$ws.Range("H62").Value = 2401.8333
$ws.Range("I62").Value = 2101.25
$ws.Range("K62").Value = 2101.25
$ws.Range("M62").Value = -1477.25
$ws.Range("H65").Value = 2401.8333
$ws.Range("I65").Value = 2101.25
$ws.Range("K65").Value = 10506.25
$ws.Range("M65").Value = -7386.25
$ws.Range("H70").Value = 3758.2
$ws.Range("J70").Value = 3697.75
$ws.Range("L70").Value = 11093.25
$ws.Range("N70").Value = -11633.25
$ws.Range("H73").Value = 3758.2
$ws.Range("J73").Value = 3697.75
$ws.Range("L73").Value = 11093.25
$ws.Range("N73").Value = -12965.25
$ws.Range("H76").Value = 3637.3635
$ws.Range("I76").Value = 3584.2083
$ws.Range("K76").Value = 3584.2083
$ws.Range("M76").Value = -3269.2083
$ws.Range("H79").Value = 3637.3635
$ws.Range("I79").Value = 3584.2083
$ws.Range("K79").Value = 3584.2083
$ws.Range("M79").Value = -2492.2083
$ws.Range("H88").Value = 2453.8572
$ws.Range("J88").Value = 1700.6666
$ws.Range("L88").Value = 1700.6666
$ws.Range("N88").Value = -2512.6666
$ws.Range("H91").Value = 2453.8572
$ws.Range("J91").Value = 1700.6666
$ws.Range("L91").Value = 1700.6666
$ws.Range("N91").Value = -4508.6666
$ws.Range("H113").Value = 4833.3335
$ws.Range("H127").Value = 2135.7
$ws.Range("I127").Value = 2238.2632
$ws.Range("K127").Value = 6714.7896
$ws.Range("M127").Value = -1754.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16450
$ws.Range("H132").Value = 717438
$ws.Range("I132").Value = 836529.6
$ws.Range("K132").Value = 2509588.8
$ws.Range("M132").Value = -2507058.8
$ws.Range("H134").Value = 64992
$ws.Range("J134").Value = 64992
$ws.Range("L134").Value = 64992
$ws.Range("N134").Value = -75132

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1844.4584
$ws.Range("J86").Value = 1965.6666
$ws.Range("L86").Value = 1965.6666
$ws.Range("N86").Value = -4211.6666
$ws.Range("H89").Value = 1844.4584
$ws.Range("J89").Value = 1965.6666
$ws.Range("L89").Value = 9828.333000000001
$ws.Range("N89").Value = -21060.333
$ws.Range("H94").Value = 1328.1818
$ws.Range("I94").Value = 1365.4
$ws.Range("K94").Value = 1365.4
$ws.Range("M94").Value = -914.4000000000001
$ws.Range("H134").Value = 5821.2812
$ws.Range("I134").Value = 2699.3447
$ws.Range("J134").Value = 36000
$ws.Range("K134").Value = 8098.034100000001
$ws.Range("L134").Value = 108000
$ws.Range("M134").Value = -5563.034100000001
$ws.Range("N134").Value = -113070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 963.1667
$ws.Range("I22").Value = 436.91428
$ws.Range("J22").Value = 2380
$ws.Range("K22").Value = 436.91428
$ws.Range("L22").Value = 2380
$ws.Range("M22").Value = -86.91428000000002
$ws.Range("N22").Value = -3080
$ws.Range("H29").Value = 700
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 1000
$ws.Range("N29").Value = -1586
$ws.Range("H94").Value = 5455.5835
$ws.Range("I94").Value = 10386.728
$ws.Range("K94").Value = 10386.728
$ws.Range("M94").Value = -9935.727999999999
$ws.Range("H107").Value = 590.5789
$ws.Range("I107").Value = 546.5
$ws.Range("K107").Value = 546.5
$ws.Range("M107").Value = 1373.5
$ws.Range("H141").Value = 512750
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 512750
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 512750
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -523110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 769.5714
$ws.Range("J132").Value = 1149.8
$ws.Range("L132").Value = 10348.2
$ws.Range("N132").Value = -15408.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 151714
$ws.Range("J24").Value = 10333
$ws.Range("L24").Value = 10333
$ws.Range("N24").Value = -10679
$ws.Range("H97").Value = 730.8049
$ws.Range("J97").Value = 167.5
$ws.Range("L97").Value = 167.5
$ws.Range("N97").Value = -1159.5
$ws.Range("H98").Value = 10642
$ws.Range("J98").Value = 10642
$ws.Range("L98").Value = 10642
$ws.Range("N98").Value = -16632
$ws.Range("H102").Value = 2027.0435
$ws.Range("I102").Value = 1936.7273
$ws.Range("K102").Value = 1936.7273
$ws.Range("M102").Value = -314.7273
$ws.Range("H109").Value = 80000
$ws.Range("J109").Value = 80000
$ws.Range("L109").Value = 80000
$ws.Range("N109").Value = -82080
$ws.Range("H132").Value = 10505.563
$ws.Range("I132").Value = 12320.643
$ws.Range("K132").Value = 36961.929
$ws.Range("M132").Value = -34431.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 372.14285
$ws.Range("I20").Value = 16.23077
$ws.Range("J20").Value = 4999
$ws.Range("K20").Value = 16.23077
$ws.Range("L20").Value = 4999
$ws.Range("M20").Value = 209.76923
$ws.Range("N20").Value = -5451
$ws.Range("H22").Value = 1744.2222
$ws.Range("J22").Value = 1744.2222
$ws.Range("L22").Value = 1744.2222
$ws.Range("N22").Value = -2334.2222
$ws.Range("H27").Value = 1744.2222
$ws.Range("J27").Value = 1744.2222
$ws.Range("L27").Value = 1744.2222
$ws.Range("N27").Value = -1958.2222
$ws.Range("H40").Value = 2676.2144
$ws.Range("I40").Value = 2042.909
$ws.Range("K40").Value = 2042.909
$ws.Range("M40").Value = -1906.909
$ws.Range("H46").Value = 4985.2
$ws.Range("I46").Value = 1963.3334
$ws.Range("J46").Value = 5740.6665
$ws.Range("K46").Value = 1963.3334
$ws.Range("L46").Value = 5740.6665
$ws.Range("M46").Value = -1775.3334
$ws.Range("N46").Value = -6116.6665
$ws.Range("H51").Value = 5250
$ws.Range("I51").Value = 5250
$ws.Range("K51").Value = 5250
$ws.Range("M51").Value = -4772
$ws.Range("H82").Value = 1361.0834
$ws.Range("I82").Value = 1750.8572
$ws.Range("J82").Value = 815.4
$ws.Range("K82").Value = 1750.8572
$ws.Range("L82").Value = 815.4
$ws.Range("M82").Value = -1389.8572
$ws.Range("N82").Value = -1537.4
$ws.Range("H85").Value = 1361.0834
$ws.Range("I85").Value = 1750.8572
$ws.Range("J85").Value = 815.4
$ws.Range("K85").Value = 1750.8572
$ws.Range("L85").Value = 815.4
$ws.Range("M85").Value = -502.8571999999999
$ws.Range("N85").Value = -3311.4
$ws.Range("H93").Value = 4043.875
$ws.Range("I93").Value = 1482.8
$ws.Range("K93").Value = 1482.8
$ws.Range("M93").Value = -234.8
$ws.Range("H132").Value = 2399.5862
$ws.Range("I132").Value = 2356.3044
$ws.Range("K132").Value = 7068.9132
$ws.Range("M132").Value = -4538.9132
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -858
$ws.Range("H21").Value = 1500000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 2101.5833
$ws.Range("J23").Value = 6666.3335
$ws.Range("L23").Value = 6666.3335
$ws.Range("N23").Value = -7124.3335
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H28").Value = 19475
$ws.Range("J28").Value = 19475
$ws.Range("L28").Value = 19475
$ws.Range("N28").Value = -20171
$ws.Range("H35").Value = 1500000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H81").Value = 2211.1365
$ws.Range("J81").Value = 2605.625
$ws.Range("L81").Value = 5211.25
$ws.Range("N81").Value = -7333.25
$ws.Range("H84").Value = 2211.1365
$ws.Range("J84").Value = 2605.625
$ws.Range("L84").Value = 26056.25
$ws.Range("N84").Value = -36664.25
$ws.Range("H112").Value = 24995.666
$ws.Range("J112").Value = 24995.666
$ws.Range("L112").Value = 24995.666
$ws.Range("N112").Value = -27949.666
$ws.Range("H132").Value = 2999.2593
$ws.Range("I132").Value = 2608.4707
$ws.Range("K132").Value = 7825.4121
$ws.Range("M132").Value = -5295.4121
